# Update "想去人数" (interested-count) figures in the F column to the
# freshly scraped values, mirroring the same event rows across the
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) tabs.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet 1) ---
$wsExhibit.Range("F3").Value  = 8532
$wsExhibit.Range("F7").Value  = 824
$wsExhibit.Range("F8").Value  = 659
$wsExhibit.Range("F9").Value  = 124
$wsExhibit.Range("F11").Value = 383
$wsExhibit.Range("F13").Value = 3684
$wsExhibit.Range("F14").Value = 263
$wsExhibit.Range("F15").Value = 147
$wsExhibit.Range("F16").Value = 791
$wsExhibit.Range("F17").Value = 774
$wsExhibit.Range("F22").Value = 1327
$wsExhibit.Range("F24").Value = 444
$wsExhibit.Range("F27").Value = 159
$wsExhibit.Range("F33").Value = 663
$wsExhibit.Range("F39").Value = 162

# --- 演出 (sheet 2) ---
$wsShow.Range("F5").Value = 2

# --- 全部类型 (sheet 4) ---
$wsAll.Range("F4").Value  = 8532
$wsAll.Range("F8").Value  = 824
$wsAll.Range("F9").Value  = 659
$wsAll.Range("F10").Value = 124
$wsAll.Range("F12").Value = 383
$wsAll.Range("F15").Value = 3684
$wsAll.Range("F16").Value = 263
$wsAll.Range("F17").Value = 147
$wsAll.Range("F19").Value = 791
$wsAll.Range("F20").Value = 774
$wsAll.Range("F25").Value = 2
$wsAll.Range("F27").Value = 1327
$wsAll.Range("F29").Value = 444
$wsAll.Range("F32").Value = 159
$wsAll.Range("F39").Value = 663
$wsAll.Range("F45").Value = 162
